$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style (bold, bordered, centered) from existing header cell H1
# into the new header cells I1 and J1, then set their text values.
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("H1").Copy($ws.Range("J1"))

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New numeric data for columns I (I0) and J (IF), rows 2-19
$dataI = @(7, 1, 3, 1, 3, 2, 1, 3, 6, 5, 3, 1, 1, 1, 1, 3, 7, 2)
$dataJ = @(9, 6, 8, 5, 6, 6, 5, 8, 8, 8, 6, 4, 4, 4, 3, 5, 7, 2)

for ($i = 0; $i -lt 18; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $dataI[$i]
    $ws.Cells.Item($row, 10).Value = $dataJ[$i]
}
